$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.077913364554649409
$ws.Range("B1").Value = 0.077913363512860337

$ws.Range("A2").Value = 0.070532356741148911
$ws.Range("B2").Value = -0.070532357843360774

$ws.Range("A3").Value = -0.084858353277130164
$ws.Range("B3").Value = 0.084858352228724759
